$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich text flattened to plain text; content matches) ---
$ws.Range("A8").Value = "Volume 29   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/21/2022  Through  11/27/2022"

# --- Row 30 cell-type changes (numeric <-> text) ---
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("D30").Value = 2
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E30").Value = -100

# --- Simple numeric value updates (Murder..Hate Crimes table, rows 14-29 + F/G/H..L30) ---
$ws.Range("F14").Value = 4
$ws.Range("H14").Value = 300
$ws.Range("I14").Value = 26
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = -3.703703703703
$ws.Range("M14").Value = -7.142857142857
$ws.Range("N14").Value = -76.785714285714
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = -83.333333333333
$ws.Range("F15").Value = 16
$ws.Range("G15").Value = 18
$ws.Range("H15").Value = -11.111111111111
$ws.Range("I15").Value = 189
$ws.Range("J15").Value = 155
$ws.Range("K15").Value = 21.935483870967
$ws.Range("L15").Value = 35.971223021582
$ws.Range("M15").Value = 52.419354838709
$ws.Range("N15").Value = 1.069518716577
$ws.Range("C16").Value = 48
$ws.Range("D16").Value = 20
$ws.Range("E16").Value = 140
$ws.Range("F16").Value = 190
$ws.Range("G16").Value = 112
$ws.Range("H16").Value = 69.642857142857
$ws.Range("I16").Value = 1724
$ws.Range("J16").Value = 1077
$ws.Range("K16").Value = 60.074280408542
$ws.Range("L16").Value = 52.296819787985
$ws.Range("M16").Value = -0.862564692351
$ws.Range("N16").Value = -78.420327950932
$ws.Range("C17").Value = 41
$ws.Range("D17").Value = 37
$ws.Range("E17").Value = 10.810810810810
$ws.Range("F17").Value = 198
$ws.Range("G17").Value = 170
$ws.Range("H17").Value = 16.470588235294
$ws.Range("I17").Value = 2270
$ws.Range("J17").Value = 1836
$ws.Range("K17").Value = 23.638344226579
$ws.Range("L17").Value = 31.670533642691
$ws.Range("M17").Value = 60.992907801418
$ws.Range("N17").Value = -14.661654135338
$ws.Range("C18").Value = 46
$ws.Range("D18").Value = 59
$ws.Range("E18").Value = -22.033898305084
$ws.Range("F18").Value = 171
$ws.Range("G18").Value = 208
$ws.Range("H18").Value = -17.788461538461
$ws.Range("I18").Value = 1807
$ws.Range("J18").Value = 1528
$ws.Range("K18").Value = 18.259162303664
$ws.Range("L18").Value = 5.119255381035
$ws.Range("M18").Value = -24.582637729549
$ws.Range("N18").Value = -86.60290628707
$ws.Range("C19").Value = 116
$ws.Range("D19").Value = 188
$ws.Range("E19").Value = -38.297872340425
$ws.Range("F19").Value = 505
$ws.Range("G19").Value = 539
$ws.Range("H19").Value = -6.307977736549
$ws.Range("I19").Value = 6389
$ws.Range("J19").Value = 4075
$ws.Range("K19").Value = 56.785276073619
$ws.Range("L19").Value = 63.820512820512
$ws.Range("M19").Value = 74.324693042292
$ws.Range("N19").Value = -16.110819327731
$ws.Range("C20").Value = 41
$ws.Range("D20").Value = 25
$ws.Range("E20").Value = 64
$ws.Range("F20").Value = 174
$ws.Range("G20").Value = 119
$ws.Range("H20").Value = 46.218487394958
$ws.Range("I20").Value = 1768
$ws.Range("J20").Value = 1262
$ws.Range("K20").Value = 40.095087163233
$ws.Range("L20").Value = 53.605560382276
$ws.Range("M20").Value = 9.950248756218
$ws.Range("N20").Value = -91.804950403263
$ws.Range("C21").Value = 293
$ws.Range("D21").Value = 335
$ws.Range("E21").Value = -12.537313432835
$ws.Range("F21").Value = 1258
$ws.Range("G21").Value = 1167
$ws.Range("H21").Value = 7.797772065124
$ws.Range("I21").Value = 14173
$ws.Range("J21").Value = 9959
$ws.Range("K21").Value = 42.313485289687
$ws.Range("L21").Value = 44.740604575163
$ws.Range("M21").Value = 29.197812215132
$ws.Range("N21").Value = -73.570656025062
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 150
$ws.Range("F22").Value = 26
$ws.Range("G22").Value = 16
$ws.Range("H22").Value = 62.5
$ws.Range("I22").Value = 280
$ws.Range("J22").Value = 143
$ws.Range("K22").Value = 95.804195804195
$ws.Range("L22").Value = 152.252252252252
$ws.Range("M22").Value = 60
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 150
$ws.Range("G23").Value = 17
$ws.Range("H23").Value = 35.294117647058
$ws.Range("I23").Value = 219
$ws.Range("J23").Value = 205
$ws.Range("K23").Value = 6.829268292682
$ws.Range("L23").Value = 14.659685863874
$ws.Range("M23").Value = 38.607594936708
$ws.Range("C24").Value = 242
$ws.Range("D24").Value = 291
$ws.Range("E24").Value = -16.838487972508
$ws.Range("F24").Value = 1194
$ws.Range("G24").Value = 1105
$ws.Range("H24").Value = 8.054298642533
$ws.Range("I24").Value = 13889
$ws.Range("J24").Value = 11057
$ws.Range("K24").Value = 25.612734014651
$ws.Range("L24").Value = 48.180945268323
$ws.Range("M24").Value = 61.125290023201
$ws.Range("C25").Value = 88
$ws.Range("D25").Value = 85
$ws.Range("E25").Value = 3.529411764705
$ws.Range("F25").Value = 401
$ws.Range("G25").Value = 368
$ws.Range("H25").Value = 8.967391304347
$ws.Range("I25").Value = 4552
$ws.Range("J25").Value = 3999
$ws.Range("K25").Value = 13.828457114278
$ws.Range("L25").Value = 31.446722494946
$ws.Range("M25").Value = 4.571559843785
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -87.5
$ws.Range("F26").Value = 22
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = -21.428571428571
$ws.Range("I26").Value = 268
$ws.Range("J26").Value = 254
$ws.Range("K26").Value = 5.511811023622
$ws.Range("L26").Value = 16.017316017316
$ws.Range("C27").Value = 9
$ws.Range("E27").Value = -30.769230769230
$ws.Range("F27").Value = 61
$ws.Range("G27").Value = 47
$ws.Range("H27").Value = 29.787234042553
$ws.Range("I27").Value = 592
$ws.Range("J27").Value = 511
$ws.Range("K27").Value = 15.851272015655
$ws.Range("L27").Value = 29.824561403508
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 12
$ws.Range("H28").Value = -25
$ws.Range("I28").Value = 72
$ws.Range("J28").Value = 76
$ws.Range("K28").Value = -5.263157894736
$ws.Range("L28").Value = -5.263157894736
$ws.Range("M28").Value = 67.441860465116
$ws.Range("N28").Value = -69.491525423728
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 8
$ws.Range("G29").Value = 8
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 61
$ws.Range("J29").Value = 54
$ws.Range("K29").Value = 12.962962962963
$ws.Range("L29").Value = 1.666666666666
$ws.Range("M29").Value = 74.285714285714
$ws.Range("N29").Value = -71.090047393364
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = -25
$ws.Range("I30").Value = 57
$ws.Range("J30").Value = 69
$ws.Range("K30").Value = -17.391304347826
$ws.Range("L30").Value = 137.5
